# Apply cryptos.xlsx update (Fri Aug  2 13:53:52 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.033.62"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3
$ws.Range("D3").Value = "3.152.70"
$ws.Range("E3").Value = "  -0.80%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.09"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.17"
$ws.Range("E6").Value = "  -1.54%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.577"
$ws.Range("E8").Value = "  -4.91%  "

# Row 9
$ws.Range("D9").Value = "3.174.12"
$ws.Range("E9").Value = "  -0.61%  "

# Row 10
$ws.Range("E10").Value = "  -0.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.66"
$ws.Range("E11").Value = "  -2.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"

# Row 13
$ws.Range("D13").Value = "3.704.89"

# Row 14
$ws.Range("E14").Value = "  -1.95%  "

# Row 15
$ws.Range("D15").Value = "65.052.33"
$ws.Range("E15").Value = "  +0.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.25"
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000157"
$ws.Range("E17").Value = "  -0.94%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.146.14"
$ws.Range("E18").Value = "  -1.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "414.60"
$ws.Range("E19").Value = "  -0.83%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -2.76%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.27"
$ws.Range("E21").Value = "  -1.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.10"
$ws.Range("E22").Value = "  -1.07%  "

# Row 23
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.07"
$ws.Range("E24").Value = "  -1.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.488"
$ws.Range("E25").Value = "  -2.41%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.196"
$ws.Range("E26").Value = "  -4.28%  "

# Row 27
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.23"
$ws.Range("E28").Value = "  +4.20%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("E30").Value = "  -0.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.49"
$ws.Range("E32").Value = "  -1.49%  "

# Row 33
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.99"
$ws.Range("E33").Value = "  -2.03%  "

# Row 34
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "163.15"
$ws.Range("E34").Value = "  +4.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.32"
$ws.Range("E35").Value = "  -1.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  +1.00%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.38"
$ws.Range("E37").Value = "  +0.70%  "

# Row 38
$ws.Range("E38").Value = "  -1.22%  "

# Row 39
$ws.Range("D39").Value = "2.618.84"
$ws.Range("E39").Value = "  -3.18%  "

# Row 40
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.11"
$ws.Range("E40").Value = "  -0.87%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.18"
$ws.Range("E41").Value = "  -1.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.45"
$ws.Range("E42").Value = "  -1.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.699"
$ws.Range("E43").Value = "  -2.68%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0624"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.39"
$ws.Range("E45").Value = "  -4.30%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.64"
$ws.Range("E46").Value = "  +0.08%  "

# Row 47
$ws.Range("E47").Value = "  -2.02%  "

# Row 48
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "292.53"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.996"
$ws.Range("E49").Value = "  -0.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0980"
$ws.Range("E50").Value = "  -1.40%  "

# Row 51
$ws.Range("E51").Value = "  -3.29%  "

Write-Output "Updated cryptos list."
